# Apply the authored change:
#  - Sheet1!C5 value 29 -> 28
#  - Active selection moves from C4 to C5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = 28
$ws.Range("C5").Select()
